$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- New data rows 100-103 ----
# Column A first for all four new rows (matches shared-string insertion order)
$ws.Cells.Item(100, 1).Value = "LEM-240-00-22KS"
$ws.Cells.Item(101, 1).Value = "LEM-250-00-A1"
$ws.Cells.Item(102, 1).Value = "LEM-276-32-3022KH"
$ws.Cells.Item(103, 1).Value = "LEM-307-00-27KS"

# Column B next
$ws.Cells.Item(100, 2).Value = "LEM-240-00-22KS"
$ws.Cells.Item(101, 2).Value = "LED-250-C00-A1"
$ws.Cells.Item(102, 2).Value = "LED-276-H70-3022"
$ws.Cells.Item(103, 2).Value = "LED-307-S00-27"
$ws.Cells.Item(101, 2).NumberFormat = "#,##0.0000"
$ws.Cells.Item(102, 2).NumberFormat = "#,##0.0000"
$ws.Cells.Item(103, 2).NumberFormat = "#,##0.0000"

# Columns C, D, E
$ws.Cells.Item(100, 3).Value = 15.6813
$ws.Cells.Item(100, 4).Value = 15.6813
$ws.Cells.Item(100, 5).Value = 40

$ws.Cells.Item(101, 3).Value = 33.4809
$ws.Cells.Item(101, 4).Value = 31.5996
$ws.Cells.Item(101, 5).Value = 80

$ws.Cells.Item(102, 3).Value = 29.1353
$ws.Cells.Item(102, 4).Value = 27.2
$ws.Cells.Item(102, 5).Value = 64

$ws.Cells.Item(103, 3).Value = 4.1327
$ws.Cells.Item(103, 4).Value = 2.8356
$ws.Cells.Item(103, 5).Value = 20

$ws.Range("C100:D103").NumberFormat = "#,##0.0000"

# ---- View: freeze top row, scroll, selection ----
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("D103").Select()
